# Insert a new weekly record row at row 300 in the "Ajo" (garlic) price
# sheet. This pushes the existing rows 300-374 down to 301-375 and adds
# one new row of data (new dimension becomes A1:R375).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 300, shifting rows 300:374
# down to 301:375 (keeping all of their original data/formatting intact).
$ws.Rows("300:300").Insert()

# Populate the newly inserted row 300 with the new weekly record.
$ws.Range("A300").Value = 5
$ws.Range("B300").Value = "Macroferia Regional de Talca"
$ws.Range("C300").Value = "Maule"
$ws.Range("D300").Value = 44855
$ws.Range("E300").Value = 7
$ws.Range("F300").Value = 100112003
$ws.Range("G300").Value = "Ajo"
$ws.Range("H300").Value = "Chino"
$ws.Range("I300").Value = "Primera"
$ws.Range("J300").Value = 300
$ws.Range("K300").Value = 19000
$ws.Range("L300").Value = 19000
$ws.Range("M300").Value = 19000
$ws.Range("N300").Value = "$/malla 10 kilos"
$ws.Range("O300").Value = "China"
$ws.Range("P300").Value = 1900
$ws.Range("Q300").Value = 10
$ws.Range("R300").Value = "Hortaliza"
